$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data row (weekly update) is inserted at row 180; all existing rows
# from 180-219 shift down by one (to 181-220). The newly inserted row 180
# receives the latest week's data, matching the constant columns shared by
# every record in this block (Terminal Hortofrutícola Agro Chillán /
# Ñuble / Zanahoria / etc.).

$ws.Rows.Item(180).Insert()

$ws.Range("A180").Value = 7
$ws.Range("B180").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C180").Value = "Ñuble"
$ws.Range("D180").Value = 44543
$ws.Range("E180").Value = 16
$ws.Range("F180").Value = 100114013
$ws.Range("G180").Value = "Zanahoria"
$ws.Range("H180").Value = "Sin especificar"
$ws.Range("I180").Value = "Primera"
$ws.Range("J180").Value = 200
$ws.Range("K180").Value = 8000
$ws.Range("L180").Value = 8500
$ws.Range("M180").Value = 8250
$ws.Range("N180").Value = "$/saco 20 kilos"
$ws.Range("O180").Value = "Provincia de Diguillín"
$ws.Range("P180").Value = 412
$ws.Range("Q180").Value = 20
$ws.Range("R180").Value = "Hortaliza"
